$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "D2"; Value = '67.790.16'; Text = $false },
    @{ Cell = "E2"; Value = '  +0.74%  '; Text = $false },
    @{ Cell = "D3"; Value = '3.816.25'; Text = $false },
    @{ Cell = "E3"; Value = '  +1.76%  '; Text = $false },
    @{ Cell = "E4"; Value = '  -0.04%  '; Text = $false },
    @{ Cell = "D5"; Value = '600.79'; Text = $true },
    @{ Cell = "E5"; Value = '  +1.25%  '; Text = $false },
    @{ Cell = "D6"; Value = '166.31'; Text = $true },
    @{ Cell = "E6"; Value = '  +0.56%  '; Text = $false },
    @{ Cell = "E8"; Value = '  +0.32%  '; Text = $false },
    @{ Cell = "E9"; Value = '  +0.99%  '; Text = $false },
    @{ Cell = "D10"; Value = '6.30'; Text = $true },
    @{ Cell = "E10"; Value = '  -0.59%  '; Text = $false },
    @{ Cell = "E11"; Value = '  +0.96%  '; Text = $false },
    @{ Cell = "E12"; Value = '  -0.15%  '; Text = $false },
    @{ Cell = "D13"; Value = '35.79'; Text = $true },
    @{ Cell = "E13"; Value = '  -0.70%  '; Text = $false },
    @{ Cell = "D14"; Value = '4.459.32'; Text = $false },
    @{ Cell = "E14"; Value = '  +1.44%  '; Text = $false },
    @{ Cell = "D15"; Value = '3.762.19'; Text = $false },
    @{ Cell = "E15"; Value = '  +0.04%  '; Text = $false },
    @{ Cell = "D16"; Value = '67.809.94'; Text = $false },
    @{ Cell = "E16"; Value = '  +0.69%  '; Text = $false },
    @{ Cell = "E17"; Value = '  +0.11%  '; Text = $false },
    @{ Cell = "D18"; Value = '7.07'; Text = $true },
    @{ Cell = "E18"; Value = '  +1.57%  '; Text = $false },
    @{ Cell = "E19"; Value = '  +0.65%  '; Text = $false },
    @{ Cell = "D20"; Value = '462.08'; Text = $true },
    @{ Cell = "E20"; Value = '  +1.76%  '; Text = $false },
    @{ Cell = "D21"; Value = '9.88'; Text = $true },
    @{ Cell = "E21"; Value = '  -0.78%  '; Text = $false },
    @{ Cell = "E22"; Value = '  +0.86%  '; Text = $false },
    @{ Cell = "E23"; Value = '  -3.14%  '; Text = $false },
    @{ Cell = "D24"; Value = '83.42'; Text = $true },
    @{ Cell = "E24"; Value = '  +0.38%  '; Text = $false },
    @{ Cell = "D25"; Value = '12.07'; Text = $true },
    @{ Cell = "E25"; Value = '  +1.89%  '; Text = $false },
    @{ Cell = "E26"; Value = '  -0.95%  '; Text = $false },
    @{ Cell = "E27"; Value = '  -0.18%  '; Text = $false },
    @{ Cell = "E28"; Value = '  +0.00%  '; Text = $false },
    @{ Cell = "D29"; Value = '3.967.75'; Text = $false },
    @{ Cell = "E29"; Value = '  +1.53%  '; Text = $false },
    @{ Cell = "E30"; Value = '  +0.59%  '; Text = $false },
    @{ Cell = "D31"; Value = '7.39'; Text = $true },
    @{ Cell = "E31"; Value = '  +2.05%  '; Text = $false },
    @{ Cell = "D32"; Value = '2.22'; Text = $true },
    @{ Cell = "E32"; Value = '  +2.62%  '; Text = $false },
    @{ Cell = "D33"; Value = '29.59'; Text = $true },
    @{ Cell = "E33"; Value = '  +0.10%  '; Text = $false },
    @{ Cell = "E34"; Value = '  +0.05%  '; Text = $false },
    @{ Cell = "B35"; Value = 'RenzoRestakedETH'; Text = $false },
    @{ Cell = "C35"; Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'; Text = $false },
    @{ Cell = "D35"; Value = '3.762.49'; Text = $false },
    @{ Cell = "E35"; Value = '  +1.52%  '; Text = $false },
    @{ Cell = "B36"; Value = 'Aptos'; Text = $false },
    @{ Cell = "C36"; Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; Text = $false },
    @{ Cell = "D36"; Value = '9.08'; Text = $true },
    @{ Cell = "E36"; Value = '  -0.85%  '; Text = $false },
    @{ Cell = "D37"; Value = '0.0999'; Text = $true },
    @{ Cell = "E37"; Value = '  +0.13%  '; Text = $false },
    @{ Cell = "E38"; Value = '  -0.69%  '; Text = $false },
    @{ Cell = "E39"; Value = '  +0.16%  '; Text = $false },
    @{ Cell = "D40"; Value = '0.998'; Text = $true },
    @{ Cell = "E40"; Value = '  +0.40%  '; Text = $false },
    @{ Cell = "E41"; Value = '  +1.28%  '; Text = $false },
    @{ Cell = "D42"; Value = '1.00'; Text = $true },
    @{ Cell = "E42"; Value = '  -0.14%  '; Text = $false },
    @{ Cell = "D44"; Value = '48.09'; Text = $true },
    @{ Cell = "E44"; Value = '  +2.45%  '; Text = $false },
    @{ Cell = "D45"; Value = '28.55'; Text = $true },
    @{ Cell = "E45"; Value = '  +10.05%  '; Text = $false },
    @{ Cell = "B46"; Value = 'TheGraph'; Text = $false },
    @{ Cell = "C46"; Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'; Text = $false },
    @{ Cell = "D46"; Value = '0.300'; Text = $true },
    @{ Cell = "E46"; Value = '  +0.88%  '; Text = $false },
    @{ Cell = "B47"; Value = 'Arweave'; Text = $false },
    @{ Cell = "C47"; Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'; Text = $false },
    @{ Cell = "D47"; Value = '43.39'; Text = $true },
    @{ Cell = "E47"; Value = '  -3.96%  '; Text = $false },
    @{ Cell = "E48"; Value = '  +12.24%  '; Text = $false },
    @{ Cell = "B49"; Value = 'Monero'; Text = $false },
    @{ Cell = "C49"; Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; Text = $false },
    @{ Cell = "D49"; Value = '148.75'; Text = $true },
    @{ Cell = "E49"; Value = '  +0.16%  '; Text = $false },
    @{ Cell = "B50"; Value = 'Cosmos'; Text = $false },
    @{ Cell = "C50"; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; Text = $false },
    @{ Cell = "D50"; Value = '8.34'; Text = $true },
    @{ Cell = "E50"; Value = '  +0.40%  '; Text = $false },
    @{ Cell = "E51"; Value = '  +0.37%  '; Text = $false }
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    if ($change.Text) {
        # Force text storage so Excel does not coerce these numeric-looking
        # strings (e.g. "600.79", "1.00", "0.300") into floating point numbers,
        # which would lose trailing zeros / precision.
        $cell.NumberFormat = "@"
        $cell.Value = $change.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $change.Value
    }
}
